# Commit: "Documents now being requested by _id rather than date"
# - Input sheet: set the Event Name and fill in the (previously blank) Date of
#   Function with an actual date value (kept as text so the underlying cell
#   stays a shared string rather than turning into a date serial number).
# - Point sheet: the placeholder "oz per shot" text values left over in the
#   Issued/Returned columns for the first few Vodka rows (and the blank
#   Non-Standard row) are replaced with real numeric 0 quantities, matching
#   every other row on the sheet.

$wb = $excel.ActiveWorkbook

$input = $wb.Worksheets.Item("Input")

$input.Range("C5").Value = "Day After Halloween"

# Keep C6 textual ("10/29/2018") instead of letting Excel coerce it into a
# date serial number: temporarily switch the number format to Text, assign
# the value, then restore the original (date) display format.
$origFormat = $input.Range("C6").NumberFormat
$input.Range("C6").NumberFormat = "@"
$input.Range("C6").Value = "10/29/2018"
$input.Range("C6").NumberFormat = $origFormat

$point = $wb.Worksheets.Item("Point")

$point.Range("E6").Value = 0
$point.Range("F6").Value = 0

$point.Range("E7").Value = 0
$point.Range("F7").Value = 0

$point.Range("E8").Value = 0
$point.Range("F8").Value = 0

$point.Range("E61").Value = 0
$point.Range("F61").Value = 0
